# Actualización automática 2025-06-09 16:40:08
# Updates static (non-formula) sales/compliance figures for
# GUERRERO FAREZ FABIAN MAURICIO across the three sheets.

$wb = $excel.ActiveWorkbook

$wsVentasGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual  = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasGrupo.Range("D18").Value = 1460.74
$wsVentasGrupo.Range("C39").Value = 518.4
$wsVentasGrupo.Range("K39").Value = 855.36

# Row 53 counters ("n de 51") shift from C/D up one column, K gains the count
$wsVentasGrupo.Range("C53").Value = "2 de 51"
$wsVentasGrupo.Range("D53").Value = "2 de 51"
$wsVentasGrupo.Range("K53").Value = "1 de 51"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual.Range("F18").Value = 1460.74
$wsVentaMensual.Range("F39").Value = 1392.23
$wsVentaMensual.Range("F53").Value = 10271.48

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento.Range("D2").Value = 1016.06
$wsCumplimiento.Range("E2").Value = 8954.283045179151
$wsCumplimiento.Range("F2").Value = 0.1019082287736613

$wsCumplimiento.Range("D3").Value = 2274.82
$wsCumplimiento.Range("E3").Value = 25182.1876
$wsCumplimiento.Range("F3").Value = 0.08285025204276085

$wsCumplimiento.Range("D15").Value = 855.36
$wsCumplimiento.Range("E15").Value = 12644.64
$wsCumplimiento.Range("F15").Value = 0.06336

$wsCumplimiento.Range("D19").Value = 10211.35
$wsCumplimiento.Range("E19").Value = 84236.09064517914
$wsCumplimiento.Range("F19").Value = 0.108116746523202
